$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price column (D) cells to Text format first so numeric-looking strings
# (e.g. "6.000", "0.07800") are preserved exactly as text, matching the source data.
$dCells = @("D2","D3","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) values
$ws.Range('D2').Value = '26.156.95'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.655.30'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '0.5298'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.2618'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = '0.06319'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').Value = '20.39'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').Value = '0.07800'
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').Value = '4.513'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '1.646.73'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '1.884.12'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '0.5488'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '0.0₅8151'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '65.39'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '26.143.62'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').Value = '4.596'
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('D21').Value = '190.69'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').Value = '10.08'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('D23').Value = '6.000'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  +4.77%  '
$ws.Range('D26').Value = '0.1225'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').Value = '7.204'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').Value = '15.96'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = '1.474'
$ws.Range('E29').Value = '  +5.29%  '
$ws.Range('D30').Value = '0.05701'
$ws.Range('E30').Value = '  -3.69%  '
$ws.Range('D31').Value = '1.274'
$ws.Range('D32').Value = '3.550'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').Value = '3.261'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('D34').Value = '1.594'
$ws.Range('E34').Value = '  +4.33%  '
$ws.Range('D35').Value = '2.805'
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').Value = '0.9491'
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('D37').Value = '2.420'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').Value = '0.5725'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').Value = '0.8511'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').Value = '5.799'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('D42').Value = '1.005'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = '104.12'
$ws.Range('E43').Value = '  +3.59%  '
$ws.Range('D44').Value = '1.036.23'
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('D45').Value = '1.797.26'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').Value = '56.72'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '0.0₈106'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').Value = '7.850'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').Value = '0.05153'
$ws.Range('E51').Value = '  +0.06%  '

# Restore default cell style on the Price column so no stray number-format
# style reference is left on the cells (keeps styling identical to original).
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
